# implement pulp instead of gurobi
# Update the ecosystem (E) and system (F) columns with newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 149.2866390832741),
    @(3, 140.2214236083316),
    @(4, 142.435132491852),
    @(5, 147.5798908804973),
    @(6, 151.4571285114751),
    @(7, 154.5689375646052),
    @(8, 159.0765148714655),
    @(9, 163.5840921783257),
    @(10, 168.091669485186),
    @(11, 172.5992467920462)
)

foreach ($row in $values) {
    $r = $row[0]
    $v = $row[1]
    $ws.Range("E$r").Value = $v
    $ws.Range("F$r").Value = $v
}
